# Hjemme passive updated meanEMG legmaxROM
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (Subj IDs) - columns B:E
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) - columns B:E
$ws.Range("B2").Value = 61.600287104867796
$ws.Range("C2").Value = 47.597757901346469
$ws.Range("D2").Value = 61.810427516813675
$ws.Range("E2").Value = 50.557737295124724

# Row 3 (STR) - columns B:E
$ws.Range("B3").Value = 62.642462350993355
$ws.Range("C3").Value = 44.981597518379402
$ws.Range("D3").Value = 66.918869894826315
$ws.Range("E3").Value = 45.291664428454979

# Match the author's updated selection
$ws.Range("B1:E3").Select()
